$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New country / file name rows (5-20), entered in the same cell order the
# author used (country names first for most rows, file names filled in
# afterwards, with a couple of later corrections).
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 1).Value  = "Algeria"
$ws.Cells.Item(5, 2).Value  = "API_DZA_DS2_en_csv_v2_5753918.csv"
$ws.Cells.Item(6, 1).Value  = "Australia"
$ws.Cells.Item(6, 2).Value  = "API_AUS_DS2_en_csv_v2_5731277.csv"

$ws.Cells.Item(8, 1).Value  = "Ghana"
$ws.Cells.Item(9, 1).Value  = "Indonesia"
$ws.Cells.Item(10, 1).Value = "Kazakhstan"

$ws.Cells.Item(13, 1).Value = "Nigeria"
$ws.Cells.Item(14, 1).Value = "Tunisia"
$ws.Cells.Item(15, 1).Value = "Turkmenistan"

$ws.Cells.Item(19, 1).Value = "Cote d'Ivoire"
$ws.Cells.Item(20, 1).Value = "Pakistan"

$ws.Cells.Item(7, 2).Value  = "API_COG_DS2_en_csv_v2_5731397.csv"
$ws.Cells.Item(8, 2).Value  = "API_GHA_DS2_en_csv_v2_5702285.csv"
$ws.Cells.Item(9, 2).Value  = "API_IDN_DS2_en_csv_v2_5707985.csv"

$ws.Cells.Item(7, 1).Value  = "Congo, Rep."

$ws.Cells.Item(11, 1).Value = "Libya"
$ws.Cells.Item(12, 1).Value = "Mexico"

$ws.Cells.Item(16, 1).Value = "United Kingdom"
$ws.Cells.Item(17, 1).Value = "United States"
$ws.Cells.Item(18, 1).Value = "Venezuela, RB"

$ws.Cells.Item(10, 2).Value = "API_KAZ_DS2_en_csv_v2_5731402.csv"
$ws.Cells.Item(11, 2).Value = "API_LBY_DS2_en_csv_v2_5731400.csv"
$ws.Cells.Item(12, 2).Value = "API_MEX_DS2_en_csv_v2_5795796.csv"
$ws.Cells.Item(13, 2).Value = "API_NGA_DS2_en_csv_v2_5693925.csv"
$ws.Cells.Item(14, 2).Value = "API_TUN_DS2_en_csv_v2_5731302.csv"
$ws.Cells.Item(15, 2).Value = "API_TKM_DS2_en_csv_v2_5695641.csv"
$ws.Cells.Item(16, 2).Value = "API_GBR_DS2_en_csv_v2_5731410.csv"
$ws.Cells.Item(17, 2).Value = "API_USA_DS2_en_csv_v2_5730402.csv"
$ws.Cells.Item(18, 2).Value = "API_VEN_DS2_en_csv_v2_5740096.csv"
$ws.Cells.Item(19, 2).Value = "API_CIV_DS2_en_csv_v2_5758464.csv"
$ws.Cells.Item(20, 2).Value = "API_PAK_DS2_en_csv_v2_5708202.csv"

# ---------------------------------------------------------------------------
# Column C ("download date") - same date as row 4 (Mozambique), 8/18/2023.
# Copy the formatting from C4 so no new number format / style gets created.
# ---------------------------------------------------------------------------
$ws.Range("C4").Copy()
for ($r = 5; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = 45156
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Expand the "country_file_name" table to cover the newly added rows
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C20"))

# ---------------------------------------------------------------------------
# Column A is now wider (longer country names) - widen it like the author did
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.5

# ---------------------------------------------------------------------------
# Leave the selection where the author's last edit left it
# ---------------------------------------------------------------------------
$ws.Range("C24").Select()
